# Append two new feed-log rows (131 and 132, 1-indexed as rows 132/133)
# to Sheet1, matching the data appended in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 132: run_id=131, rss_url_id=1, date=2024-06-17 17:12:06, response=200, item_count=22
$ws.Range("A132").Value = 131
$ws.Range("B132").Value = 1
$ws.Range("C132").Value = "2024-06-17 17:12:06"
$ws.Range("D132").Value = 200
$ws.Range("E132").Value = 22

# Row 133: run_id=132, rss_url_id=2, date=2024-06-17 17:12:06, response=200, item_count=3
$ws.Range("A133").Value = 132
$ws.Range("B133").Value = 2
$ws.Range("C133").Value = "2024-06-17 17:12:06"
$ws.Range("D133").Value = 200
$ws.Range("E133").Value = 3
